$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend new rows 17-19 by copying formatting (style) from row 16 (A and B columns use style s="1")
$ws.Range("A16").Copy($ws.Range("A17:A19"))

$ws.Range("A10").Value2 = 8
$ws.Range("B10").Value2 = "Gaussian-Quadrature"
$ws.Range("C10").Value2 = 1.061031675377856
$ws.Range("D10").Value2 = 1.233246932143259
$ws.Range("E10").Value2 = 1.000829744036772
$ws.Range("F10").Value2 = 1.061031675377856
$ws.Range("G10").Value2 = 0.8675645298150684
$ws.Range("H10").Value2 = 1.451569371944928
$ws.Range("I10").Value2 = 0.962324404813809
$ws.Range("J10").Value2 = 1.233246932143259
$ws.Range("K10").Value2 = 1.117038338090016
$ws.Range("L10").Value2 = 1.089035006733936
$ws.Range("M10").Value2 = 1.096094443021949

$ws.Range("A11").Value2 = 9
$ws.Range("B11").Value2 = "Spiral-90deg-10rot-5space"
$ws.Range("C11").Value2 = 0.6835093538392826
$ws.Range("D11").Value2 = 1.270775534408722
$ws.Range("E11").Value2 = 1.15623494300233
$ws.Range("F11").Value2 = 0.6835093538392826
$ws.Range("G11").Value2 = 0.8487452942885764
$ws.Range("H11").Value2 = 1.848710910479194
$ws.Range("I11").Value2 = 0.9912362659532078
$ws.Range("J11").Value2 = 1.270775534408722
$ws.Range("K11").Value2 = 1.213505238705526
$ws.Range("L11").Value2 = 0.9485072962724043
$ws.Range("M11").Value2 = 1.133202050328552

$ws.Range("A12").Value2 = 10
$ws.Range("B12").Value2 = "Spiral-90deg-15rot-5space"
$ws.Range("C12").Value2 = 0.6844387052542044
$ws.Range("D12").Value2 = 1.25634367426999
$ws.Range("E12").Value2 = 1.159231955148341
$ws.Range("F12").Value2 = 0.6844387052542044
$ws.Range("G12").Value2 = 0.8438655564722771
$ws.Range("H12").Value2 = 1.85457730012915
$ws.Range("I12").Value2 = 0.9938762640849956
$ws.Range("J12").Value2 = 1.25634367426999
$ws.Range("K12").Value2 = 1.207787814709165
$ws.Range("L12").Value2 = 0.9461132599816849
$ws.Range("M12").Value2 = 1.13205557589316

$ws.Range("A13").Value2 = 11
$ws.Range("B13").Value2 = "Spiral-90deg-10rot-3space"
$ws.Range("C13").Value2 = 0.6828290851119831
$ws.Range("D13").Value2 = 1.268264768497686
$ws.Range("E13").Value2 = 1.15716135833466
$ws.Range("F13").Value2 = 0.6828290851119831
$ws.Range("G13").Value2 = 0.8475327214295615
$ws.Range("H13").Value2 = 1.850873208256092
$ws.Range("I13").Value2 = 0.9920919898878994
$ws.Range("J13").Value2 = 1.268264768497686
$ws.Range("K13").Value2 = 1.212713063416173
$ws.Range("L13").Value2 = 0.9477710742640779
$ws.Range("M13").Value2 = 1.133125521919647

$ws.Range("A14").Value2 = 12
$ws.Range("B14").Value2 = "NoRotation-tilt60deg"
$ws.Range("C14").Value2 = 2.283984000000002
$ws.Range("D14").Value2 = 0.02368399999999932
$ws.Range("E14").Value2 = 0.9084639999999992
$ws.Range("F14").Value2 = 2.283984000000002
$ws.Range("G14").Value2 = 0.3111519999999989
$ws.Range("H14").Value2 = 1.408540000000001
$ws.Range("I14").Value2 = 1.113191999999998
$ws.Range("J14").Value2 = 0.02368399999999932
$ws.Range("K14").Value2 = 0.4660739999999993
$ws.Range("L14").Value2 = 1.375029000000001
$ws.Range("M14").Value2 = 1.008169333333333

$ws.Range("A15").Value2 = 13
$ws.Range("B15").Value2 = "Rotation-NoTilt"
$ws.Range("C15").Value2 = 3.435012499999988
$ws.Range("D15").Value2 = 0.01
$ws.Range("E15").Value2 = 0.45
$ws.Range("F15").Value2 = 3.435012499999988
$ws.Range("G15").Value2 = 0.3886250000000002
$ws.Range("H15").Value2 = 0.22
$ws.Range("I15").Value2 = 0.9836250000000027
$ws.Range("J15").Value2 = 0.01
$ws.Range("K15").Value2 = 0.23
$ws.Range("L15").Value2 = 1.832506249999994
$ws.Range("M15").Value2 = 0.9145437499999985

$ws.Range("A16").Value2 = 14
$ws.Range("B16").Value2 = "Rotation-60detTilt"
$ws.Range("C16").Value2 = 2.385723970457597
$ws.Range("D16").Value2 = 0.4234165041152013
$ws.Range("E16").Value2 = 0.6843187093503985
$ws.Range("F16").Value2 = 2.385723970457597
$ws.Range("G16").Value2 = 0.6454906476544017
$ws.Range("H16").Value2 = 0.558515681280005
$ws.Range("I16").Value2 = 0.997835114086404
$ws.Range("J16").Value2 = 0.4234165041152013
$ws.Range("K16").Value2 = 0.5538676067327999
$ws.Range("L16").Value2 = 1.469795788595198
$ws.Range("M16").Value2 = 0.9492167711573344

$ws.Range("A17").Value2 = 15
$ws.Range("B17").Value2 = "HexGrid-90degTilt5degRes"
$ws.Range("C17").Value2 = 0.9862841412949093
$ws.Range("D17").Value2 = 0.9987428813913282
$ws.Range("E17").Value2 = 0.990868888415937
$ws.Range("F17").Value2 = 0.9862841412949093
$ws.Range("G17").Value2 = 0.9907493041183262
$ws.Range("H17").Value2 = 1.001139772509748
$ws.Range("I17").Value2 = 0.9922728475483893
$ws.Range("J17").Value2 = 0.9987428813913282
$ws.Range("K17").Value2 = 0.9948058849036325
$ws.Range("L17").Value2 = 0.990545013099271
$ws.Range("M17").Value2 = 0.9933429725464397

$ws.Range("A18").Value2 = 16
$ws.Range("B18").Value2 = "HexGrid-90degTilt22p5degRes"
$ws.Range("C18").Value2 = 1.068201957836057
$ws.Range("D18").Value2 = 0.9139794098522952
$ws.Range("E18").Value2 = 0.9943948781468208
$ws.Range("F18").Value2 = 1.068201957836057
$ws.Range("G18").Value2 = 0.9903156575817584
$ws.Range("H18").Value2 = 0.841388336655014
$ws.Range("I18").Value2 = 0.9980923855235098
$ws.Range("J18").Value2 = 0.9139794098522952
$ws.Range("K18").Value2 = 0.954187143999558
$ws.Range("L18").Value2 = 1.011194550917807
$ws.Range("M18").Value2 = 0.9677287709325758

$ws.Range("A19").Value2 = 17
$ws.Range("B19").Value2 = "HexGrid-60degTilt5degRes"
$ws.Range("C19").Value2 = 0.9941038032747677
$ws.Range("D19").Value2 = 1.140738590472542
$ws.Range("E19").Value2 = 0.954947285193842
$ws.Range("F19").Value2 = 0.9941038032747677
$ws.Range("G19").Value2 = 1.040374737777881
$ws.Range("H19").Value2 = 0.9737487568172937
$ws.Range("I19").Value2 = 0.9640482619635485
$ws.Range("J19").Value2 = 1.140738590472542
$ws.Range("K19").Value2 = 1.047842937833192
$ws.Range("L19").Value2 = 1.02097337055398
$ws.Range("M19").Value2 = 1.011326905916646

Write-Output "done"